$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.22
$ws.Range("R2").Value = 1.67
$ws.Range("U2").Value = 5.2
$ws.Range("V2").Value = 1.16

# Row 4
$ws.Range("G4").Value = 3.1
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 3.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.78
$ws.Range("S4").Value = 2.88
$ws.Range("T4").Value = 1.4
$ws.Range("U4").Value = 4.8
$ws.Range("V4").Value = 1.19
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 1.13
$ws.Range("AA4").Value = 2.25
$ws.Range("AB4").Value = 1.57
$ws.Range("AC4").Value = 6.5
$ws.Range("AG4").Value = 34
$ws.Range("AI4").Value = 6
$ws.Range("AP4").Value = 26

# Row 5
$ws.Range("U5").Value = 4.7
$ws.Range("V5").Value = 1.19

# Row 6
$ws.Range("G6").Value = 1.4
$ws.Range("I6").Value = 8
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 1.47
$ws.Range("R6").Value = 2.7
$ws.Range("S6").Value = 1.95
$ws.Range("T6").Value = 1.9
$ws.Range("U6").Value = 2.65
$ws.Range("V6").Value = 1.49
$ws.Range("Y6").Value = 1.4
$ws.Range("Z6").Value = 2.75
$ws.Range("AA6").Value = 2.2
$ws.Range("AB6").Value = 1.62
$ws.Range("AH6").Value = 34
$ws.Range("AI6").Value = 10
$ws.Range("AM6").Value = 17
$ws.Range("AO6").Value = 23
$ws.Range("AP6").Value = 101

# Row 8
$ws.Range("G8").Value = 1.95
$ws.Range("I8").Value = 3.4
$ws.Range("K8").Value = 2.3
$ws.Range("L8").Value = 3.75
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.33
$ws.Range("S8").Value = 1.7
$ws.Range("T8").Value = 2.1
$ws.Range("W8").Value = 2.63
$ws.Range("X8").Value = 1.44
$ws.Range("Y8").Value = 1.33
$ws.Range("Z8").Value = 3.25
$ws.Range("AA8").Value = 1.57
$ws.Range("AB8").Value = 2.25
$ws.Range("AC8").Value = 9.5
$ws.Range("AD8").Value = 11
$ws.Range("AF8").Value = 19
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 13
$ws.Range("AK8").Value = 12
$ws.Range("AQ8").Value = 26
$ws.Range("AR8").Value = 29
$ws.Range("AS8").Value = 126

# Row 9
$ws.Range("G9").Value = 1.48
$ws.Range("H9").Value = 4.33
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2.5
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 5
$ws.Range("S9").Value = 1.57
$ws.Range("T9").Value = 2.35
$ws.Range("U9").Value = 1.93
$ws.Range("V9").Value = 1.88
$ws.Range("W9").Value = 2.38
$ws.Range("X9").Value = 1.53
$ws.Range("Y9").Value = 1.29
$ws.Range("Z9").Value = 3.5
$ws.Range("AC9").Value = 9
$ws.Range("AJ9").Value = 8.5

# Row 11
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("O11").Value = 1.22
$ws.Range("P11").Value = 4
$ws.Range("S11").Value = 1.7
$ws.Range("T11").Value = 2.1
$ws.Range("W11").Value = 2.75
$ws.Range("X11").Value = 1.4

# Row 12
$ws.Range("G12").Value = 3.5
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 2.15
$ws.Range("L12").Value = 2.88
$ws.Range("N12").Value = 9.5
$ws.Range("T12").Value = 1.75
$ws.Range("Y12").Value = 1.44
$ws.Range("Z12").Value = 2.63
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 17
$ws.Range("AI12").Value = 9.5
$ws.Range("AR12").Value = 29

# Row 14
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("S14").Value = 1.88
$ws.Range("T14").Value = 1.98
$ws.Range("W14").Value = 3.25
$ws.Range("X14").Value = 1.33

# Row 15
$ws.Range("Q15").Value = 1.75
$ws.Range("R15").Value = 2.05
$ws.Range("Y15").Value = 1.53
$ws.Range("Z15").Value = 2.38
$ws.Range("AA15").Value = 2.1
$ws.Range("AB15").Value = 1.67

# Row 16
$ws.Range("G16").Value = 2.9
$ws.Range("I16").Value = 2.6
$ws.Range("J16").Value = 3.75
$ws.Range("L16").Value = 3.4
$ws.Range("AC16").Value = 7.5
$ws.Range("AD16").Value = 13
$ws.Range("AE16").Value = 12
$ws.Range("AG16").Value = 29
$ws.Range("AN16").Value = 11
$ws.Range("AP16").Value = 26
$ws.Range("AQ16").Value = 23

# Row 17
$ws.Range("G17").Value = 2.55
$ws.Range("H17").Value = 3.55
$ws.Range("I17").Value = 2.47
$ws.Range("J17").Value = 2.95
$ws.Range("K17").Value = 2.35
$ws.Range("L17").Value = 2.9
$ws.Range("N17").Value = 9.5
$ws.Range("T17").Value = 2.65
$ws.Range("W17").Value = 2
$ws.Range("X17").Value = 1.72
$ws.Range("Z17").Value = 3.55
$ws.Range("AB17").Value = 2.87
$ws.Range("AC17").Value = 15
$ws.Range("AD17").Value = 18.5
$ws.Range("AF17").Value = 32
$ws.Range("AH17").Value = 18
$ws.Range("AI17").Value = 9.5
$ws.Range("AJ17").Value = 7.9
$ws.Range("AL17").Value = 27
$ws.Range("AM17").Value = 15
$ws.Range("AN17").Value = 18
$ws.Range("AR17").Value = 17.5

# Row 18
$ws.Range("G18").Value = 4.4
$ws.Range("I18").Value = 1.7
$ws.Range("J18").Value = 4.4
$ws.Range("K18").Value = 2.35
$ws.Range("L18").Value = 2.18
$ws.Range("O18").Value = 1.17
$ws.Range("P18").Value = 4.45
$ws.Range("S18").Value = 1.53
$ws.Range("T18").Value = 2.35
$ws.Range("W18").Value = 2.25
$ws.Range("X18").Value = 1.57
$ws.Range("Y18").Value = 1.29
$ws.Range("Z18").Value = 3.3
$ws.Range("AA18").Value = 1.55
$ws.Range("AB18").Value = 2.32
$ws.Range("AC18").Value = 17.5
$ws.Range("AD18").Value = 30
$ws.Range("AE18").Value = 14
$ws.Range("AF18").Value = 75
$ws.Range("AG18").Value = 35
$ws.Range("AH18").Value = 32
$ws.Range("AJ18").Value = 8
$ws.Range("AM18").Value = 9.75
$ws.Range("AN18").Value = 9.75
$ws.Range("AO18").Value = 8
$ws.Range("AP18").Value = 14
$ws.Range("AQ18").Value = 11.75
$ws.Range("AR18").Value = 18.5

# Row 19
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.9
$ws.Range("J19").Value = 3.4
$ws.Range("K19").Value = 1.95
$ws.Range("L19").Value = 3.75
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 7
$ws.Range("AD19").Value = 11
$ws.Range("AF19").Value = 26
$ws.Range("AM19").Value = 7.5
$ws.Range("AN19").Value = 13
$ws.Range("AO19").Value = 12
$ws.Range("AQ19").Value = 29
